# remove-products-expiration-date.py ready for verification
#
# Append a new product row (barcode, count, expiration date) below the
# existing data, reusing the formatting of the row above it, then move
# the active selection the way it was left in the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the number formatting (barcode style on col A, date style on col C)
# down from row 3 into the new row 4 before writing the new values.
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)  # xlPasteFormats

# New product entry.
$ws.Range("A4").Value = 7610700949085
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 44593

# Leave the selection where the author left it.
$null = $ws.Range("D7").Select()
